# Update countries & provincias Spain
# - Update the "last updated" timestamp in A1
# - Update India's (row 14) case numbers
# - Update Nepal's case numbers; since Nepal's total cases now exceed
#   those of "Sudan del Sur" and "Congo", the table (sorted descending
#   by total cases) re-sorts, shifting those two rows down by one and
#   moving Nepal up to row 131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp update
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 06:05"

# India - row 14
$ws.Range("B14").Value = 118501
$ws.Range("C14").Value = 275
$ws.Range("E14").Value = 66363
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 3585

# Rows 131-133 now hold (after re-sort by total cases desc):
#   131 = Nepal (updated numbers, moved up)
#   132 = Sudan del Sur (previously row 131, unchanged numbers)
#   133 = Congo (previously row 132, unchanged numbers)

$ws.Range("A131").Value = "Nepal"
$ws.Range("B131").Value = 487
$ws.Range("C131").Value = 30
$ws.Range("D131").Value = 49
$ws.Range("E131").Value = 435
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 3

$ws.Range("A132").Value = "Sudan del Sur"
$ws.Range("B132").Value = 481
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 4
$ws.Range("E132").Value = 473
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 4

$ws.Range("A133").Value = "Congo"
$ws.Range("B133").Value = 469
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 137
$ws.Range("E133").Value = 316
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 16
